# Trade #79 closed at 2026-02-17 08:58:44 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.5   # Current Capital
$wsSummary.Range("B4").Value = 0.51     # Total P&L $
$wsSummary.Range("B6").Value = 79       # Total Trades
$wsSummary.Range("B7").Value = 33       # Winning Trades
$wsSummary.Range("B9").Value = 41.77    # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.5     # Capital
$wsStatus.Range("D4").Value = 79        # Trades
$wsStatus.Range("E4").Value = 0.51      # P&L $
$wsStatus.Range("F4").Value = 0.5       # P&L %
$wsStatus.Range("G4").Value = 41.77     # Win Rate %

# --- Append new trade row (row 80) to "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(80, 1).Value = 79

    # Date/Time columns are stored as plain text in this workbook, not Excel
    # dates - force text format before assignment, then clear the leftover
    # format stamp so no extra style is introduced.
    $ws.Cells.Item(80, 2).NumberFormat = "@"
    $ws.Cells.Item(80, 2).Value = "2026-02-17"
    $ws.Cells.Item(80, 2).ClearFormats()

    $ws.Cells.Item(80, 3).Value = "08:58:38"
    $ws.Cells.Item(80, 4).Value = "MarketMaking"
    $ws.Cells.Item(80, 5).Value = "UP"
    $ws.Cells.Item(80, 6).Value = 0.1
    $ws.Cells.Item(80, 7).Value = 0.11
    $ws.Cells.Item(80, 8).Value = "CLOSED"
    $ws.Cells.Item(80, 9).Value = 10
    $ws.Cells.Item(80, 10).Value = 0.01
    $ws.Cells.Item(80, 11).Value = 100.5
    $ws.Cells.Item(80, 12).Value = 0
    $ws.Cells.Item(80, 13).Value = 0
    $ws.Cells.Item(80, 14).Value = 0.6
    $ws.Cells.Item(80, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(80, 16).Value = "early_exit"
    $ws.Cells.Item(80, 17).Value = 0.14
}

Write-Output "Applied trade #79 update"
